$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.546.03"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.815.34"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  -0.72%  "
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "306.44"
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3594"
$ws.Range("E8").Value = "  -2.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.37"
$ws.Range("E9").Value = "  +2.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07097"
$ws.Range("E10").Value = "  -0.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8938"
$ws.Range("E11").Value = "  +1.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07767"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.30"
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("D14").Value = "1.833.13"
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.260"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.314"
$ws.Range("E16").Value = "  -1.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "85.25"
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008602"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("D21").Value = "26.577.47"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.14"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.960"
$ws.Range("E23").Value = "  -1.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.50"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.967"
$ws.Range("E25").Value = "  -0.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.71"
$ws.Range("E26").Value = "  +0.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.81"
$ws.Range("E27").Value = "  -1.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.040"
$ws.Range("E28").Value = "  -2.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.44"
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.841"
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08720"
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.126"
$ws.Range("E32").Value = "  +2.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7505"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.743"
$ws.Range("E34").Value = "  +6.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.424"
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.112"
$ws.Range("E36").Value = "  -0.68%  "
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.910"
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05100"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5110"
$ws.Range("E41").Value = "  +1.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.772"
$ws.Range("E42").Value = "  -2.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1511"
$ws.Range("E43").Value = "  -3.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.035"
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4706"
$ws.Range("E45").Value = "  +1.57%  "
$ws.Range("E46").Value = "  -0.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.953"
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "99.86"
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.572"
$ws.Range("E49").Value = "  -1.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05983"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.96"
$ws.Range("E51").Value = "  -0.47%  "
